$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the table header / timestamp column name (C1) to the new value.
$ws.Range("C1").Value = "2024-11-15 02:23:39"

# Swap the status values for two rows.
$ws.Range("C6").Value = "Puntual"
$ws.Range("C9").Value = "Retardo"
